$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize "de"/"del"/"los"/"la" particles in specific municipality/state names
$ws.Range("B12").Value = "Villa De Álvarez"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A23").Value = "Estado De México"
$ws.Range("B23").Value = "Tlalnepantla De Baz"
$ws.Range("B33").Value = "Acapulco De Juárez"
$ws.Range("B35").Value = "Atlamajalcingo Del Monte"
$ws.Range("B38").Value = "Chilpancingo De Los Bravo"
$ws.Range("B40").Value = "Cutzamala De Pinzón"
$ws.Range("B41").Value = "Zihuatanejo De Azueta"
$ws.Range("B45").Value = "Tlapa De Comonfort"
$ws.Range("B49").Value = "Progreso De Obregón"
$ws.Range("B52").Value = "Lagos De Moreno"
$ws.Range("B53").Value = "Tamazula De Gordiano"
$ws.Range("B57").Value = "Unión De Tula"
$ws.Range("B78").Value = "Santo Domingo De Morelos"
$ws.Range("B86").Value = "Amealco De Bonfil"
$ws.Range("B87").Value = "Cadereyta De Montes"
$ws.Range("B89").Value = "Landa De Matamoros"
$ws.Range("B93").Value = "Axtla De Terrazas"
$ws.Range("B110").Value = "Martínez De La Torre"

# Remove trailing footer/metadata rows that are no longer part of the table
# (delete bottom-most block first so row numbers of the upper block stay valid)
$ws.Range("A476:D480").EntireRow.Delete()
$ws.Range("A124:D128").EntireRow.Delete()

# Update used range dimension to match new data extent
$ws.UsedRange | Out-Null
